# BE_LTD_Dispatcher/Data/SummaryReport.xlsx
#
# The dispatcher process now drops each day's bank statements into a
# dated sub-folder (e.g. "06-23-2022") under ...\Data\Input\, instead of
# directly under ...\Data\Input\<Bank Name>\. Update every recorded
# "Attachment Path" value in the summary report to reflect the new,
# dated folder layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldToNew = @{
    "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\SBI Bank\Vidya Sagar Reddy.pdf"   = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\06-23-2022\SBI Bank\Vidya Sagar Reddy.pdf";
    "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\Axis Bank\Vidya Sagar Reddy.pdf"  = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\06-23-2022\Axis Bank\Vidya Sagar Reddy.pdf";
    "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\ICICI Bank\Vidya Sagar Reddy.pdf" = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\06-23-2022\ICICI Bank\Vidya Sagar Reddy.pdf";
    "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\HDFC Bank\Vidya Sagar Reddy.pdf"  = "C:\Users\TEMP\OneDrive - bradsol.com\Documents\GitHub\brad-uipath\BE_LTD_Dispatcher\Data\Input\06-23-2022\HDFC Bank\Vidya Sagar Reddy.pdf";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($val -ne $null -and $oldToNew.ContainsKey($val)) {
        $cell.Value = $oldToNew[$val]
    }
}

# The "Attachment Path" column is wider now that the paths contain the
# extra "06-23-2022\" segment, and the "Month" column picked up its own
# explicit (auto-fit) width too.
$ws.Columns.Item(3).ColumnWidth = 127.16666666666667   # -> stored width 128
$ws.Columns.Item(9).ColumnWidth = 8.333333333333334    # -> stored width ~9.17 (best fit)

# The selected whole-column range moved one column to the right (I -> J).
$ws.Range("J1:J1048576").Select() | Out-Null
